# This workbook's single table is being reshaped:
#  - old layout: rows = years (2019, 2018, 2017), columns = metrics
#  - new layout: rows = metrics, columns = years (2019, 2018, 2017),
#    with a new generic header row ("Unnamed: 0..3") and a new
#    "Number of employees" metric row holding the years themselves.
# The result is effectively a brand new 4 (cols) x 6 (rows) table, so we
# clear out the old 5x4 table and write the new values/styles directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused rightmost column (E) entirely.
$ws.Columns("E").Delete()

# Wipe old contents from the area the old table occupied, and strip the old
# bold/bordered "key column" styling from A2:A4 (rows 5:6 are brand new and
# already style-less) so we start from a clean slate before laying out the
# new table. A1 deliberately keeps its existing bold/bordered style (style
# index 1 in the original file) since the new header row needs that same
# style too.
$ws.Range("A1:D4").ClearContents()
$ws.Range("A2:D6").ClearFormats()

# ---- Row 1: generic header labels ----
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"
$ws.Range("D1").Value = "Unnamed: 3"

# Re-apply the workbook's existing bold/bordered header style across the
# whole header row by copying it from A1, which already carries that style,
# onto B1:D1. Using copy/paste-special of formats (instead of re-building
# the formatting by hand) lets the engine reuse the existing style record
# instead of creating a new one.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Row 2: Number of employees ----
$ws.Range("A2").Value = "Number of employees"
$ws.Range("B2").Value = 2019
$ws.Range("C2").Value = 2018
# ---- Row 3: Total participants ----
$ws.Range("A3").Value = "Total participants"
$ws.Range("B3").Value = 784
$ws.Range("C3").Value = 991
# ---- Row 4: Total for company sports activities ----
$ws.Range("A4").Value = "Total for company sports activities"
$ws.Range("B4").Value = 554
$ws.Range("C4").Value = 683
# ---- Row 5: Total for sickness/injury prevention courses ----
$ws.Range("A5").Value = "Total for sickness/injury prevention courses"
$ws.Range("B5").Value = 230
$ws.Range("C5").Value = 308
# ---- Row 6: Average health rate for the year (%) ----
$ws.Range("A6").Value = "Average health rate for the year (%)"
$ws.Range("B6").Value = 95.3
$ws.Range("C6").Value = 95.3

# Column D holds the 2017 figures, which (as in the original workbook) are
# stored as text rather than numbers. Force text storage by pre-formatting
# the column as Text before assigning, then strip the formatting back off
# (the underlying stored values stay text; only the now-unneeded number
# format is removed) so the cells end up with the default (unstyled) look,
# matching every other plain data cell in the sheet.
$ws.Range("D2:D6").NumberFormat = "@"
$ws.Range("D2").Value = "2017"
$ws.Range("D3").Value = "1,061"
$ws.Range("D4").Value = "740"
$ws.Range("D5").Value = "321"
$ws.Range("D6").Value = "96.5"
$ws.Range("D2:D6").ClearFormats()
